$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, new D (Price) text or $null if unchanged, new E (Volume 1h) text
$updates = @(
    ,@(2, '69.759.41', '  +0.63%  ')
    ,@(3, '2.508.62', '  +0.21%  ')
    ,@(4, $null, '  +0.01%  ')
    ,@(5, '576.04', '  +0.07%  ')
    ,@(6, '167.07', '  +0.58%  ')
    ,@(7, $null, '  -0.04%  ')
    ,@(8, $null, '  -0.23%  ')
    ,@(9, '2.506.54', '  +0.21%  ')
    ,@(10, '0.163', '  +4.06%  ')
    ,@(11, '0.168', '  +0.06%  ')
    ,@(12, '0.356', '  +4.22%  ')
    ,@(13, $null, '  +2.23%  ')
    ,@(14, '2.967.34', '  +0.24%  ')
    ,@(15, '0.0000179', '  +3.03%  ')
    ,@(16, '69.673.95', '  +0.58%  ')
    ,@(17, '24.90', '  +0.81%  ')
    ,@(18, '2.502.25', '  +0.01%  ')
    ,@(19, '11.24', '  -1.29%  ')
    ,@(20, '7.51', '  -3.07%  ')
    ,@(21, '349.88', '  +0.85%  ')
    ,@(22, $null, '  -0.24%  ')
    ,@(23, $null, '  +0.91%  ')
    ,@(24, $null, '  -0.11%  ')
    ,@(25, '70.57', '  +3.49%  ')
    ,@(26, $null, '  +0.08%  ')
    ,@(27, '8.83', '  -0.59%  ')
    ,@(28, '2.639.23', '  +0.05%  ')
    ,@(29, '0.994', '  -2.16%  ')
    ,@(30, '0.0₃0895', '  +0.12%  ')
    ,@(31, $null, '  +0.27%  ')
    ,@(32, '460.25', '  -0.88%  ')
    ,@(33, $null, '  -2.68%  ')
    ,@(34, $null, '  -0.48%  ')
    ,@(35, '0.999', '  +0.05%  ')
    ,@(36, $null, '  +2.43%  ')
    ,@(37, '158.31', '  +3.56%  ')
    ,@(38, '19.07', '  +0.68%  ')
    ,@(39, '18.51', '  +0.88%  ')
    ,@(40, $null, '  +0.03%  ')
    ,@(41, $null, '  +1.41%  ')
    ,@(42, '4.70', '  -0.41%  ')
    ,@(43, $null, '  +1.41%  ')
    ,@(44, $null, '  -0.02%  ')
    ,@(45, '2.21', '  -3.70%  ')
    ,@(46, '1.08', '  -6.14%  ')
    ,@(47, '141.59', '  -0.65%  ')
    ,@(48, $null, '  -0.32%  ')
    ,@(49, $null, '  -0.80%  ')
    ,@(50, '0.0736', '  +0.77%  ')
    ,@(51, $null, '  -0.64%  ')
)

foreach ($entry in $updates) {
    $row = $entry[0]
    $dVal = $entry[1]
    $eVal = $entry[2]
    if ($null -ne $dVal) {
        $dCell = $ws.Range("D$row")
        $dCell.NumberFormat = '@'
        $dCell.Value = $dVal
        $dCell.Style = 'Normal'
    }
    $ws.Range("E$row").Value = $eVal
}
